$d = $word.ActiveDocument

# ======================================================================
# 1. "four" -> "three"  (paragraph: "...You will be working with four files...")
# ======================================================================
$full = $d.Content.Text
$i = $full.IndexOf("working with four files")
$fourStart = $i + "working with ".Length
$rngFour = $d.Range($fourStart, $fourStart + "four".Length)
$rngFour.Text = "three"

# ======================================================================
# 2. Replace "box.js, and the-all-in-one.html." with "and box.js."
#    i.e. remove the ", and the-all-in-one.html" run/mention, and turn the
#    ", " before box.js into ", and ".
# ======================================================================
$full = $d.Content.Text
$i = $full.IndexOf(", and the-all-in-one.html")
$rngRemove = $d.Range($i, $i + ", and the-all-in-one.html".Length)
$rngRemove.Text = ""

$full = $d.Content.Text
$i = $full.IndexOf("box.js")
$rngComma = $d.Range($i - 2, $i)
$rngComma.Text = ", and "

# ======================================================================
# 3. Remove the (collapsed) "_GoBack" bookmark that currently sits right
#    after "...and he" (before " has no idea ...").
# ======================================================================
$bmOld = $d.Bookmarks("_GoBack")
$bmOld.Delete()

# ======================================================================
# 4. Delete the table row documenting "the-all-in-one.html" (along with its
#    "/1" mark), and bump the "box.js" row's mark from "2" to "3".
# ======================================================================
$t = $d.Tables(1)
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $rowText = $t.Cell($r, 1).Range.Text
    if ($rowText -like "*the-all-in-one.html*") {
        $t.Rows($r).Delete()
        break
    }
}

$boxRowIndex = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $rowText = $t.Cell($r, 1).Range.Text
    if ($rowText -like "box.js*") {
        $boxRowIndex = $r
        break
    }
}
$markCell = $t.Cell($boxRowIndex, 2)
$cellStart = $markCell.Range.Start
$cellEnd = $markCell.Range.End
for ($p = $cellStart; $p -lt $cellEnd; $p++) {
    $ch = $d.Range($p, $p + 1).Text
    if ($ch -eq "2") {
        $rngMark = $d.Range($p, $p + 1)
        $rngMark.Text = "3"
        break
    }
}

# ======================================================================
# 5. Add a new (collapsed) "_GoBack" bookmark right after the "2" mark in
#    the "whatever.js is fully filled in" row (first data row's mark cell).
#    A direct zero-length Range placed immediately before a paragraph mark
#    cannot be used with Bookmarks.Add reliably, so a temporary character is
#    inserted to anchor the position, the bookmark is created next to it,
#    and the temporary character is removed again.
# ======================================================================
$whateverRowIndex = -1
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $rowText = $t.Cell($r, 1).Range.Text
    if ($rowText -like "whatever.js*") {
        $whateverRowIndex = $r
        break
    }
}
$markCell2 = $t.Cell($whateverRowIndex, 2)
$insertPos = $markCell2.Range.End - 1
$insPoint = $d.Range($insertPos, $insertPos)
$insPoint.InsertAfter("X")

$rngBm = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $rngBm)

$rngX = $d.Range($insertPos, $insertPos + 1)
$rngX.Text = ""

Write-Output "done"
